$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domains")
$ws.Columns("H").Insert()
$ws.Range("H2:H6").Clear()
$ws.Range("H8").Clear()
$ws.Range("H1").Interior.Color = 65535
$ws.Range("H7").Interior.Color = 65535
Write-Host "done"
